$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "平潭发展"
$ws.Range("B3").Value = "平潭发展"
$ws.Range("B4").Value = "雷科防务"
$ws.Range("C4").Value = "雷科防务"
$ws.Range("A5").Value = "海欣食品"
$ws.Range("B5").Value = "合富中国"
$ws.Range("A6").Value = "中兴通讯"
$ws.Range("B6").Value = "实达集团"
$ws.Range("C6").Value = "合富中国"
$ws.Range("A7").Value = "赛微电子"
$ws.Range("B7").Value = "中兴通讯"
$ws.Range("C7").Value = "海欣食品"
$ws.Range("A8").Value = "合富中国"
$ws.Range("B8").Value = "赛微电子"
$ws.Range("C8").Value = "海王生物"
$ws.Range("A9").Value = "道明光学"
$ws.Range("B9").Value = "海欣食品"
$ws.Range("C9").Value = "海峡创新"
$ws.Range("A10").Value = "实达集团"
$ws.Range("B10").Value = "福蓉科技"
$ws.Range("C10").Value = "赛微电子"
$ws.Range("A11").Value = "海峡创新"
$ws.Range("B11").Value = "道明光学"
$ws.Range("C11").Value = "海马汽车"
$ws.Range("A12").Value = "榕基软件"
$ws.Range("B12").Value = "海马汽车"
$ws.Range("A13").Value = "海王生物"
$ws.Range("B13").Value = "海王生物"
$ws.Range("C13").Value = "榕基软件"
$ws.Range("A14").Value = "海马汽车"
$ws.Range("B14").Value = "榕基软件"
$ws.Range("C14").Value = "福蓉科技"
$ws.Range("A15").Value = "广和通"
$ws.Range("B15").Value = "华夏幸福"
$ws.Range("C15").Value = "国晟科技"
$ws.Range("A16").Value = "冰山冷热"
$ws.Range("B16").Value = "通宇通讯"
$ws.Range("C16").Value = "国风新材"
$ws.Range("A17").Value = "福蓉科技"
$ws.Range("B17").Value = "招标股份"
$ws.Range("C17").Value = "金富科技"
$ws.Range("A18").Value = "东百集团"
$ws.Range("B18").Value = "海峡创新"
$ws.Range("C18").Value = "顺灏股份"
$ws.Range("A19").Value = "国风新材"
$ws.Range("B19").Value = "东百集团"
$ws.Range("C19").Value = "中兴通讯"
$ws.Range("A20").Value = "招标股份"
$ws.Range("C20").Value = "特发信息"
$ws.Range("A21").Value = "华夏幸福"
$ws.Range("B21").Value = "人民同泰"
$ws.Range("C21").Value = "东百集团"
